# COREESG_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer text (A10) from
#    2021-06-10 to 2021-06-14
#  - refresh the Weight (D) / Percent Change (E) figures for rows 2-7
#
# The sheet is protected, so we need to unprotect it before writing, then
# restore protection afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Disclaimer text (shared string used by A10) -------------------------
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."
$ws.Range("A10").Value = $disclaimer

# --- Weight / Percent Change figures --------------------------------------
$ws.Range("D2").Value = 0.2525033436932321
$ws.Range("E2").Value = 0.005790574231944623

$ws.Range("D3").Value = 0.5383848054072217
$ws.Range("E3").Value = -0.003923620193565402

$ws.Range("D4").Value = 0.05049205321327568
$ws.Range("E4").Value = 0.005147759771210625

$ws.Range("D5").Value = 0.09895593190488744
$ws.Range("E5").Value = -0.00782937365010794

$ws.Range("D6").Value = 0.05966386578138319
$ws.Range("E6").Value = -0.007197382769901783

$ws.Range("E7").Value = -0.001594543823975836

# --- Restore sheet protection ---------------------------------------------
$ws.Protect("D382")
